# Refresh market-price derived columns (H:N) across the Leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 55
$ws_ALC.Range("H55").Value = 275.33334
$ws_ALC.Range("I55").Value = 154.28572
$ws_ALC.Range("J55").Value = 444.8
$ws_ALC.Range("K55").Value = 154.28572
$ws_ALC.Range("L55").Value = 444.8
$ws_ALC.Range("M55").Value = 59.71428
$ws_ALC.Range("N55").Value = -872.8

# ALC row 98
$ws_ALC.Range("H98").Value = 1189.7826
$ws_ALC.Range("I98").Value = 1103.85
$ws_ALC.Range("J98").Value = 1762.6666
$ws_ALC.Range("K98").Value = 1103.85
$ws_ALC.Range("L98").Value = 1762.6666
$ws_ALC.Range("M98").Value = 394.1500000000001
$ws_ALC.Range("N98").Value = -4758.6666

# ALC row 122
$ws_ALC.Range("H122").Value = 1189.7826
$ws_ALC.Range("I122").Value = 1103.85
$ws_ALC.Range("J122").Value = 1762.6666
$ws_ALC.Range("K122").Value = 3311.55
$ws_ALC.Range("L122").Value = 5287.9998
$ws_ALC.Range("M122").Value = -861.5499999999997
$ws_ALC.Range("N122").Value = -10187.9998

# ALC row 137
$ws_ALC.Range("H137").Value = 1234.6666
$ws_ALC.Range("I137").Value = 500
$ws_ALC.Range("J137").Value = 1326.5
$ws_ALC.Range("K137").Value = 1500
$ws_ALC.Range("L137").Value = 3979.5
$ws_ALC.Range("M137").Value = 1050
$ws_ALC.Range("N137").Value = -9079.5

# ALC row 138
$ws_ALC.Range("H138").Value = 3236.8333
$ws_ALC.Range("I138").Value = 2764.3
$ws_ALC.Range("J138").Value = 5599.5
$ws_ALC.Range("K138").Value = 8292.900000000001
$ws_ALC.Range("L138").Value = 16798.5
$ws_ALC.Range("M138").Value = -3152.900000000001
$ws_ALC.Range("N138").Value = -27078.5

# ARM row 61
$ws_ARM.Range("H61").Value = 3632.1724
$ws_ARM.Range("I61").Value = 3014.2632
$ws_ARM.Range("J61").Value = 4806.2
$ws_ARM.Range("K61").Value = 3014.2632
$ws_ARM.Range("L61").Value = 4806.2
$ws_ARM.Range("M61").Value = -2802.2632
$ws_ARM.Range("N61").Value = -5230.2

# ARM row 136
$ws_ARM.Range("H136").Value = 3632.1724
$ws_ARM.Range("I136").Value = 3014.2632
$ws_ARM.Range("J136").Value = 4806.2
$ws_ARM.Range("K136").Value = 9042.7896
$ws_ARM.Range("L136").Value = 14418.6
$ws_ARM.Range("M136").Value = -6492.7896
$ws_ARM.Range("N136").Value = -19518.6

# CRP row 7
$ws_CRP.Range("H7").Value = 210.2
$ws_CRP.Range("I7").Value = 164.22223
$ws_CRP.Range("J7").Value = 279.16666
$ws_CRP.Range("K7").Value = 164.22223
$ws_CRP.Range("L7").Value = 279.16666
$ws_CRP.Range("M7").Value = -51.22223
$ws_CRP.Range("N7").Value = -505.16666

# CRP row 31
$ws_CRP.Range("H31").Value = 8635.596
$ws_CRP.Range("I31").Value = 2669.5
$ws_CRP.Range("J31").Value = 12338.689
$ws_CRP.Range("K31").Value = 2669.5
$ws_CRP.Range("L31").Value = 12338.689
$ws_CRP.Range("M31").Value = -2374.5
$ws_CRP.Range("N31").Value = -12928.689

# CRP row 34
$ws_CRP.Range("H34").Value = 8635.596
$ws_CRP.Range("I34").Value = 2669.5
$ws_CRP.Range("J34").Value = 12338.689
$ws_CRP.Range("K34").Value = 2669.5
$ws_CRP.Range("L34").Value = 12338.689
$ws_CRP.Range("M34").Value = -2467.5
$ws_CRP.Range("N34").Value = -12742.689

# CRP row 134
$ws_CRP.Range("H134").Value = 4037.2727
$ws_CRP.Range("I134").Value = 3962.2222
$ws_CRP.Range("J134").Value = 4375
$ws_CRP.Range("K134").Value = 11886.6666
$ws_CRP.Range("L134").Value = 13125
$ws_CRP.Range("M134").Value = -9351.6666
$ws_CRP.Range("N134").Value = -18195

# CUL row 56
$ws_CUL.Range("H56").Value = 11921.77
$ws_CUL.Range("I56").Value = 11921.77
$ws_CUL.Range("K56").Value = 11921.77
$ws_CUL.Range("M56").Value = -11391.77

# CUL row 68
$ws_CUL.Range("H68").Value = 3569
$ws_CUL.Range("I68").Value = 3096.25
$ws_CUL.Range("J68").Value = 3805.375
$ws_CUL.Range("K68").Value = 9288.75
$ws_CUL.Range("L68").Value = 11416.125
$ws_CUL.Range("M68").Value = -8477.75
$ws_CUL.Range("N68").Value = -13038.125

# CUL row 71
$ws_CUL.Range("H71").Value = 3569
$ws_CUL.Range("I71").Value = 3096.25
$ws_CUL.Range("J71").Value = 3805.375
$ws_CUL.Range("K71").Value = 27866.25
$ws_CUL.Range("L71").Value = 34248.375
$ws_CUL.Range("M71").Value = -23810.25
$ws_CUL.Range("N71").Value = -42360.375

# CUL row 92
$ws_CUL.Range("H92").Value = 1183.1111
$ws_CUL.Range("I92").Value = 117.5
$ws_CUL.Range("J92").Value = 2035.6
$ws_CUL.Range("K92").Value = 352.5
$ws_CUL.Range("L92").Value = 6106.799999999999
$ws_CUL.Range("M92").Value = 895.5
$ws_CUL.Range("N92").Value = -8602.799999999999

# GSM row 132
$ws_GSM.Range("H132").Value = 4090.724
$ws_GSM.Range("I132").Value = 3293.6155
$ws_GSM.Range("J132").Value = 10999
$ws_GSM.Range("K132").Value = 9880.8465
$ws_GSM.Range("L132").Value = 32997
$ws_GSM.Range("M132").Value = -7350.8465
$ws_GSM.Range("N132").Value = -38057

# LTW row 7
$ws_LTW.Range("H7").Value = 43482030
$ws_LTW.Range("I7").Value = 52634856
$ws_LTW.Range("J7").Value = 6115.5
$ws_LTW.Range("K7").Value = 52634856
$ws_LTW.Range("L7").Value = 6115.5
$ws_LTW.Range("M7").Value = -52634744
$ws_LTW.Range("N7").Value = -6339.5

# LTW row 22
$ws_LTW.Range("H22").Value = 3274.6924
$ws_LTW.Range("I22").Value = 1924.7142
$ws_LTW.Range("J22").Value = 4849.6665
$ws_LTW.Range("K22").Value = 1924.7142
$ws_LTW.Range("L22").Value = 4849.6665
$ws_LTW.Range("M22").Value = -1629.7142
$ws_LTW.Range("N22").Value = -5439.6665

# LTW row 27
$ws_LTW.Range("H27").Value = 3274.6924
$ws_LTW.Range("I27").Value = 1924.7142
$ws_LTW.Range("J27").Value = 4849.6665
$ws_LTW.Range("K27").Value = 1924.7142
$ws_LTW.Range("L27").Value = 4849.6665
$ws_LTW.Range("M27").Value = -1817.7142
$ws_LTW.Range("N27").Value = -5063.6665

# LTW row 126
$ws_LTW.Range("H126").Value = 43482030
$ws_LTW.Range("I126").Value = 52634856
$ws_LTW.Range("J126").Value = 6115.5
$ws_LTW.Range("K126").Value = 157904568
$ws_LTW.Range("L126").Value = 18346.5
$ws_LTW.Range("M126").Value = -157902098
$ws_LTW.Range("N126").Value = -23286.5

# LTW row 132
$ws_LTW.Range("H132").Value = 11056.789
$ws_LTW.Range("I132").Value = 10138.6
$ws_LTW.Range("J132").Value = 14500
$ws_LTW.Range("K132").Value = 30415.8
$ws_LTW.Range("L132").Value = 43500
$ws_LTW.Range("M132").Value = -27885.8
$ws_LTW.Range("N132").Value = -48560

# LTW row 139
$ws_LTW.Range("H139").Value = 99706
$ws_LTW.Range("J139").Value = 99706
$ws_LTW.Range("L139").Value = 99706
$ws_LTW.Range("N139").Value = -109986

# LTW row 141
$ws_LTW.Range("H141").Value = 0
$ws_LTW.Range("J141").Value = 0
$ws_LTW.Range("L141").Value = 0
$ws_LTW.Range("N141").ClearContents()

# WVR row 132
$ws_WVR.Range("H132").Value = 51602.5
$ws_WVR.Range("I132").Value = 3200
$ws_WVR.Range("J132").Value = 100005
$ws_WVR.Range("K132").Value = 9600
$ws_WVR.Range("L132").Value = 300015
$ws_WVR.Range("M132").Value = -7070
$ws_WVR.Range("N132").Value = -305075

# WVR row 139
$ws_WVR.Range("H139").Value = 73331.32000000001
$ws_WVR.Range("J139").Value = 73331.32000000001
$ws_WVR.Range("L139").Value = 73331.32000000001
$ws_WVR.Range("N139").Value = -83611.32000000001
